$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextValue($ws, $cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue $ws 'D2' '34.869.27'
Set-TextValue $ws 'E2' '  -2.59%  '
Set-TextValue $ws 'D3' '1.804.77'
Set-TextValue $ws 'E3' '  -3.07%  '
Set-TextValue $ws 'E4' '  +0.31%  '
Set-TextValue $ws 'D5' '230.92'
Set-TextValue $ws 'E5' '  -0.18%  '
Set-TextValue $ws 'E6' '  -1.79%  '
Set-TextValue $ws 'D7' '1.01'
Set-TextValue $ws 'E7' '  +0.32%  '
Set-TextValue $ws 'D8' '38.91'
Set-TextValue $ws 'E8' '  -8.91%  '
Set-TextValue $ws 'E9' '  +2.94%  '
Set-TextValue $ws 'E10' '  -3.14%  '
Set-TextValue $ws 'E11' '  -1.45%  '
Set-TextValue $ws 'D12' '2.065.43'
Set-TextValue $ws 'D13' '1.811.76'
Set-TextValue $ws 'E13' '  -2.87%  '
Set-TextValue $ws 'E14' '  -3.72%  '
Set-TextValue $ws 'D15' '10.84'
Set-TextValue $ws 'E15' '  -6.14%  '
Set-TextValue $ws 'E16' '  -5.54%  '
Set-TextValue $ws 'D17' '34.789.69'
Set-TextValue $ws 'E17' '  -2.66%  '
Set-TextValue $ws 'D18' '68.97'
Set-TextValue $ws 'E18' '  -2.00%  '
Set-TextValue $ws 'E19' '  -3.58%  '
Set-TextValue $ws 'D20' '238.20'
Set-TextValue $ws 'E20' '  -3.41%  '
Set-TextValue $ws 'D21' '11.70'
Set-TextValue $ws 'E21' '  -4.71%  '
Set-TextValue $ws 'E22' '  -4.34%  '
Set-TextValue $ws 'E23' '  +0.18%  '
Set-TextValue $ws 'E24' '  -0.81%  '
Set-TextValue $ws 'D25' '171.14'
Set-TextValue $ws 'E25' '  +0.47%  '
Set-TextValue $ws 'E26' '  -3.87%  '
Set-TextValue $ws 'D27' '17.20'
Set-TextValue $ws 'E27' '  -4.19%  '
Set-TextValue $ws 'E28' '  -3.92%  '
Set-TextValue $ws 'E29' '  +5.49%  '
Set-TextValue $ws 'E30' '  +0.31%  '
Set-TextValue $ws 'D31' '3.96'
Set-TextValue $ws 'E31' '  +0.27%  '
Set-TextValue $ws 'E32' '  +0.38%  '
Set-TextValue $ws 'D33' '3.90'
Set-TextValue $ws 'E33' '  -4.06%  '
Set-TextValue $ws 'E34' '  -8.15%  '
Set-TextValue $ws 'E35' '  +3.71%  '
Set-TextValue $ws 'D36' '0.680'
Set-TextValue $ws 'E36' '  -1.45%  '
Set-TextValue $ws 'D37' '90.13'
Set-TextValue $ws 'E37' '  -3.97%  '
Set-TextValue $ws 'E38' '  +1.55%  '
Set-TextValue $ws 'D39' '1.302.67'
Set-TextValue $ws 'E39' '  -3.24%  '
Set-TextValue $ws 'E40' '  -3.14%  '
Set-TextValue $ws 'E41' '  -1.03%  '
Set-TextValue $ws 'D42' '14.39'
Set-TextValue $ws 'E42' '  -5.35%  '
Set-TextValue $ws 'D43' '0.953'
Set-TextValue $ws 'E43' '  -7.53%  '
Set-TextValue $ws 'B44' 'RenderToken'
Set-TextValue $ws 'C44' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws 'D44' '2.18'
Set-TextValue $ws 'E44' '  -13.57%  '
Set-TextValue $ws 'B45' 'MXToken'
Set-TextValue $ws 'C45' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws 'D45' '2.69'
Set-TextValue $ws 'E45' '  -4.82%  '
Set-TextValue $ws 'D46' '6.11'
Set-TextValue $ws 'E46' '  -0.10%  '
Set-TextValue $ws 'E47' '  -2.43%  '
Set-TextValue $ws 'D48' '1.984.78'
Set-TextValue $ws 'E49' '  +0.26%  '
Set-TextValue $ws 'D50' '0.0670'
Set-TextValue $ws 'E50' '  +7.43%  '
Set-TextValue $ws 'D51' '98.53'
Set-TextValue $ws 'E51' '  -6.08%  '
